$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 5863462
$ws.Range("I33").Value = 8207480.5
$ws.Range("K33").Value = 8207480.5
$ws.Range("M33").Value = -8207251.5

$ws.Range("H86").Value = 95242690
$ws.Range("I86").Value = 111115576
$ws.Range("J86").Value = 66671504
$ws.Range("K86").Value = 111115576
$ws.Range("L86").Value = 66671504
$ws.Range("M86").Value = -111114453
$ws.Range("N86").Value = -66673750

$ws.Range("H89").Value = 95242690
$ws.Range("I89").Value = 111115576
$ws.Range("J89").Value = 66671504
$ws.Range("K89").Value = 555577880
$ws.Range("L89").Value = 333357520
$ws.Range("M89").Value = -555572264
$ws.Range("N89").Value = -333368752

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 194496.53
$ws.Range("I32").Value = 233465.72
$ws.Range("K32").Value = 233465.72
$ws.Range("M32").Value = -233178.72

$ws.Range("H61").Value = 1072865.9
$ws.Range("I61").Value = 3188.0244
$ws.Range("K61").Value = 3188.0244
$ws.Range("M61").Value = -2976.0244

$ws.Range("H94").Value = 9998
$ws.Range("J94").Value = 9998
$ws.Range("L94").Value = 9998
$ws.Range("N94").Value = -11800

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H96").Value = 29614.666
$ws.Range("J96").Value = 29614.666
$ws.Range("L96").Value = 29614.666
$ws.Range("N96").Value = -35106.666

$ws.Range("H136").Value = 1072865.9
$ws.Range("I136").Value = 3188.0244
$ws.Range("K136").Value = 9564.073199999999
$ws.Range("M136").Value = -7014.073199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 890.25
$ws.Range("I25").Value = 890.25
$ws.Range("K25").Value = 890.25
$ws.Range("M25").Value = -655.25

$ws.Range("H94").Value = 10496.9375
$ws.Range("I94").Value = 11995.182
$ws.Range("K94").Value = 11995.182
$ws.Range("M94").Value = -11544.182

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H134").Value = 29034740
$ws.Range("I134").Value = 2639.9
$ws.Range("K134").Value = 7919.700000000001
$ws.Range("M134").Value = -5384.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 89999.39999999999
$ws.Range("J20").Value = 89999.39999999999
$ws.Range("L20").Value = 89999.39999999999
$ws.Range("N20").Value = -90471.39999999999

$ws.Range("H30").Value = 89999.39999999999
$ws.Range("J30").Value = 89999.39999999999
$ws.Range("L30").Value = 89999.39999999999
$ws.Range("N30").Value = -90181.39999999999

$ws.Range("H31").Value = 2759.0393
$ws.Range("I31").Value = 2355.2727
$ws.Range("K31").Value = 2355.2727
$ws.Range("M31").Value = -2060.2727

$ws.Range("H34").Value = 2759.0393
$ws.Range("I34").Value = 2355.2727
$ws.Range("K34").Value = 2355.2727
$ws.Range("M34").Value = -2153.2727

$ws.Range("H43").Value = 14988.667
$ws.Range("J43").Value = 14988.667
$ws.Range("L43").Value = 14988.667
$ws.Range("N43").Value = -15356.667

$ws.Range("H101").Value = 14988.667
$ws.Range("J101").Value = 14988.667
$ws.Range("L101").Value = 14988.667
$ws.Range("N101").Value = -21478.667

$ws.Range("H106").Value = 64994.5
$ws.Range("J106").Value = 64994.5
$ws.Range("L106").Value = 64994.5
$ws.Range("N106").Value = -67518.5

$ws.Range("H128").Value = 89999.39999999999
$ws.Range("J128").Value = 89999.39999999999
$ws.Range("L128").Value = 89999.39999999999
$ws.Range("N128").Value = -99959.39999999999

$ws.Range("H132").Value = 19036.61
$ws.Range("I132").Value = 23038.086
$ws.Range("J132").Value = 3364.1667
$ws.Range("K132").Value = 69114.258
$ws.Range("L132").Value = 10092.5001
$ws.Range("M132").Value = -66584.258
$ws.Range("N132").Value = -15152.5001

$ws.Range("H134").Value = 1315.4166
$ws.Range("I134").Value = 1001.3958
$ws.Range("K134").Value = 3004.1874
$ws.Range("M134").Value = -469.1873999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 2480.7144
$ws.Range("J50").Value = 3329.2
$ws.Range("L50").Value = 9987.599999999999
$ws.Range("N50").Value = -10949.6

$ws.Range("H53").Value = 2480.7144
$ws.Range("J53").Value = 3329.2
$ws.Range("L53").Value = 9987.599999999999
$ws.Range("N53").Value = -10949.6

$ws.Range("H103").Value = 1435728.4
$ws.Range("I103").Value = 3334700
$ws.Range("J103").Value = 11499.75
$ws.Range("K103").Value = 10004100
$ws.Range("L103").Value = 34499.25
$ws.Range("M103").Value = -10003221
$ws.Range("N103").Value = -36257.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 10122.5
$ws.Range("J22").Value = 10122.5
$ws.Range("L22").Value = 10122.5
$ws.Range("N22").Value = -11180.5

$ws.Range("H41").Value = 7857.5557
$ws.Range("I41").Value = 7857.5557
$ws.Range("K41").Value = 7857.5557
$ws.Range("M41").Value = -7502.5557

$ws.Range("H80").Value = 62652500
$ws.Range("I80").Value = 186665.83
$ws.Range("K80").Value = 186665.83
$ws.Range("M80").Value = -185667.83

$ws.Range("H83").Value = 62652500
$ws.Range("I83").Value = 186665.83
$ws.Range("K83").Value = 933329.1499999999
$ws.Range("M83").Value = -928337.1499999999

$ws.Range("H97").Value = 4312.0386
$ws.Range("I97").Value = 465.9565
$ws.Range("K97").Value = 465.9565
$ws.Range("M97").Value = 30.04349999999999

$ws.Range("H122").Value = 1979.5946
$ws.Range("I122").Value = 1756.3939
$ws.Range("K122").Value = 5269.1817
$ws.Range("M122").Value = -2819.1817

$ws.Range("H123").Value = 58999.668
$ws.Range("J123").Value = 58999.668
$ws.Range("L123").Value = 58999.668
$ws.Range("N123").Value = -63899.668

$ws.Range("H133").Value = 73749.625
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 19899
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 19899
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 19899
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -20125

$ws.Range("H16").Value = 1526.7916
$ws.Range("I16").Value = 1352.15
$ws.Range("K16").Value = 1352.15
$ws.Range("M16").Value = -1182.15

$ws.Range("H22").Value = 4566.9287
$ws.Range("I22").Value = 1699
$ws.Range("J22").Value = 5925.421
$ws.Range("K22").Value = 1699
$ws.Range("L22").Value = 5925.421
$ws.Range("M22").Value = -1404
$ws.Range("N22").Value = -6515.421

$ws.Range("H27").Value = 4566.9287
$ws.Range("I27").Value = 1699
$ws.Range("J27").Value = 5925.421
$ws.Range("K27").Value = 1699
$ws.Range("L27").Value = 5925.421
$ws.Range("M27").Value = -1592
$ws.Range("N27").Value = -6139.421

$ws.Range("H28").Value = 19899
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 19899
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 19899
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -20363

$ws.Range("H37").Value = 19899
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 19899
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 19899
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -20113

$ws.Range("H46").Value = 11264.059
$ws.Range("I46").Value = 17299.285
$ws.Range("J46").Value = 7039.4
$ws.Range("K46").Value = 17299.285
$ws.Range("L46").Value = 7039.4
$ws.Range("M46").Value = -17111.285
$ws.Range("N46").Value = -7415.4

$ws.Range("H55").Value = 1659.7037
$ws.Range("I55").Value = 1826.1428
$ws.Range("K55").Value = 1826.1428
$ws.Range("M55").Value = -1653.1428

$ws.Range("H93").Value = 1323.0834
$ws.Range("I93").Value = 1388.3
$ws.Range("J93").Value = 997
$ws.Range("K93").Value = 1388.3
$ws.Range("L93").Value = 997
$ws.Range("M93").Value = -140.3
$ws.Range("N93").Value = -3493

$ws.Range("H100").Value = 2824.6667
$ws.Range("I100").Value = 2704.6365
$ws.Range("K100").Value = 2704.6365
$ws.Range("M100").Value = -2163.6365

$ws.Range("H103").Value = 28333
$ws.Range("J103").Value = 28333
$ws.Range("L103").Value = 28333
$ws.Range("N103").Value = -30677

$ws.Range("H122").Value = 3341.625
$ws.Range("I122").Value = 2914.4546
$ws.Range("K122").Value = 8743.363799999999
$ws.Range("M122").Value = -6293.363799999999

$ws.Range("H132").Value = 3146.7827
$ws.Range("I132").Value = 2899.4707
$ws.Range("K132").Value = 8698.4121
$ws.Range("M132").Value = -6168.4121

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 15534.857
$ws.Range("I61").Value = 15534.857
$ws.Range("K61").Value = 15534.857
$ws.Range("M61").Value = -15242.857

$ws.Range("H75").Value = 28251.143
$ws.Range("I75").Value = 19499.6
$ws.Range("K75").Value = 19499.6
$ws.Range("M75").Value = -18563.6

$ws.Range("H78").Value = 28251.143
$ws.Range("I78").Value = 19499.6
$ws.Range("K78").Value = 58498.8
$ws.Range("M78").Value = -53818.8

$ws.Range("H132").Value = 2165.224
$ws.Range("I132").Value = 1693.878
$ws.Range("K132").Value = 5081.634
$ws.Range("M132").Value = -2551.634

$ws.Range("H136").Value = 79317.53999999999
$ws.Range("I136").Value = 101314.4
$ws.Range("J136").Value = 5994.6665
$ws.Range("K136").Value = 303943.2
$ws.Range("L136").Value = 17983.9995
$ws.Range("M136").Value = -301393.2
$ws.Range("N136").Value = -23083.9995
